$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.741.21'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.629.68'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.58'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.34'
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('D7').Value = '3.631.19'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.96'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '4.243.74'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000210'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.99'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '3.621.05'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '66.827.75'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.62'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.14'
$ws.Range('E21').Value = '  +1.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '429.36'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.622'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('E26').Value = '  +4.13%  '
$ws.Range('E27').Value = '  +5.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.57'
$ws.Range('E28').Value = '  +5.22%  '
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '3.627.12'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.51'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.91'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.70'
$ws.Range('E37').Value = '  +1.24%  '
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '176.91'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0864'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.903'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '46.16'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.60'
$ws.Range('E45').Value = '  +8.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.18'
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.13'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.04'
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('E51').Value = '  +1.55%  '
